$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"22.68000000000011"
$ws.Range("H2").Value = [double]"1.069257884456221e-07"
$ws.Range("I2").Value = [double]"1.069257884456221e-07"
$ws.Range("L2").Value = [double]"60.69499867193912"
$ws.Range("M2").Value = "[39.89949835398596, 81.49049898989227]"
$ws.Range("N2").Value = [double]"4.728680096377502e-07"
$ws.Range("O2").Value = [double]"4.728680096377502e-07"
$ws.Range("P2").Value = [double]"1.742184514603348"
$ws.Range("Q2").Value = "[1.3648160276856558, 2.119553001521041]"
$ws.Range("R2").Value = [double]"4.825695398835705e-12"
$ws.Range("S2").Value = [double]"4.825695398835705e-12"
$ws.Range("T2").Value = [double]"70.37303827211969"
$ws.Range("U2").Value = "[58.006697841073304, 82.73937870316607]"
$ws.Range("V2").Value = [double]"5.995204332975845e-15"
$ws.Range("W2").Value = [double]"5.995204332975845e-15"
$ws.Range("X2").Value = [double]"16.39135135135143"
$ws.Range("Y2").Value = [double]"15.02918918918926"
$ws.Range("Z2").Value = [double]"17.7535135135136"
# Row 3
$ws.Range("F3").Value = [double]"22.68000000000011"
$ws.Range("H3").Value = [double]"0.0009811331950526636"
$ws.Range("I3").Value = [double]"0.0009811331950526636"
$ws.Range("L3").Value = [double]"40.69627153579664"
$ws.Range("M3").Value = "[15.05322229947997, 66.33932077211331]"
$ws.Range("N3").Value = [double]"0.002545522675897249"
$ws.Range("O3").Value = [double]"0.002545522675897249"
$ws.Range("P3").Value = [double]"1.817658211986886"
$ws.Range("Q3").Value = "[1.0880791372793457, 2.5472372866944273]"
$ws.Range("R3").Value = [double]"8.666754596209358e-06"
$ws.Range("S3").Value = [double]"8.666754596209358e-06"
$ws.Range("T3").Value = [double]"51.79590645540019"
$ws.Range("U3").Value = "[37.68692512881637, 65.904887781984]"
$ws.Range("V3").Value = [double]"2.663316012174732e-09"
$ws.Range("W3").Value = [double]"2.663316012174732e-09"
$ws.Range("X3").Value = [double]"16.118918918919"
$ws.Range("Y3").Value = [double]"13.48540540540547"
$ws.Range("Z3").Value = [double]"18.75243243243253"
# Row 4
$ws.Range("F4").Value = [double]"22.68000000000011"
$ws.Range("H4").Value = [double]"3.869065203776501e-07"
$ws.Range("I4").Value = [double]"3.869065203776501e-07"
$ws.Range("L4").Value = [double]"50.76076302849162"
$ws.Range("M4").Value = "[28.973371505503053, 72.5481545514802]"
$ws.Range("N4").Value = [double]"2.541694821278639e-05"
$ws.Range("O4").Value = [double]"2.541694821278639e-05"
$ws.Range("P4").Value = [double]"2.484342538874811"
$ws.Range("Q4").Value = "[2.0692372032653488, 2.8994478744842738]"
$ws.Range("R4").Value = [double]"1.110223024625157e-15"
$ws.Range("S4").Value = [double]"1.110223024625157e-15"
$ws.Range("T4").Value = [double]"69.91328244904889"
$ws.Range("U4").Value = "[58.58675067828315, 81.23981421981463]"
$ws.Range("V4").Value = [double]"4.440892098500626e-16"
$ws.Range("W4").Value = [double]"4.440892098500626e-16"
$ws.Range("X4").Value = [double]"13.7124324324325"
$ws.Range("Y4").Value = [double]"12.21405405405411"
$ws.Range("Z4").Value = [double]"15.21081081081088"
# Row 5
$ws.Range("F5").Value = [double]"22.92000000000014"
$ws.Range("H5").Value = [double]"0.0005310030143512723"
$ws.Range("I5").Value = [double]"0.0005310030143512723"
$ws.Range("L5").Value = [double]"38.01962813957645"
$ws.Range("M5").Value = "[15.752726988051386, 60.28652929110151]"
$ws.Range("N5").Value = [double]"0.001269305596920844"
$ws.Range("O5").Value = [double]"0.001269305596920844"
$ws.Range("P5").Value = [double]"3.025237370123504"
$ws.Range("Q5").Value = "[2.283079345852042, 3.767395394394967]"
$ws.Range("R5").Value = [double]"1.708573282854786e-10"
$ws.Range("S5").Value = [double]"1.708573282854786e-10"
$ws.Range("T5").Value = [double]"59.99617800338493"
$ws.Range("U5").Value = "[46.927099428817996, 73.06525657795186]"
$ws.Range("V5").Value = [double]"5.71120928327673e-12"
$ws.Range("W5").Value = [double]"5.71120928327673e-12"
$ws.Range("X5").Value = [double]"11.88444444444452"
$ws.Range("Y5").Value = [double]"9.177177177177235"
$ws.Range("Z5").Value = [double]"14.5917117117118"
# Row 6
$ws.Range("F6").Value = [double]"22.92000000000014"
$ws.Range("H6").Value = [double]"0.00139465062563926"
$ws.Range("I6").Value = [double]"0.00139465062563926"
$ws.Range("L6").Value = [double]"38.37748350688009"
$ws.Range("M6").Value = "[12.884855670112884, 63.870111343647295]"
$ws.Range("N6").Value = [double]"0.004020451474406883"
$ws.Range("O6").Value = [double]"0.004020451474406883"
$ws.Range("P6").Value = [double]"2.849132076228582"
$ws.Range("Q6").Value = "[2.0943951023931957, 3.6038690500639676]"
$ws.Range("R6").Value = [double]"1.310888952943401e-09"
$ws.Range("S6").Value = [double]"1.310888952943401e-09"
$ws.Range("T6").Value = [double]"58.88937810968183"
$ws.Range("U6").Value = "[44.64449706192819, 73.13425915743547]"
$ws.Range("V6").Value = [double]"1.15991216631528e-10"
$ws.Range("W6").Value = [double]"1.15991216631528e-10"
$ws.Range("X6").Value = [double]"12.52684684684692"
$ws.Range("Y6").Value = [double]"9.773693693693751"
$ws.Range("Z6").Value = [double]"15.28000000000009"
# Row 7
$ws.Range("F7").Value = [double]"22.92000000000014"
$ws.Range("H7").Value = [double]"0.0007090700484622126"
$ws.Range("I7").Value = [double]"0.0007090700484622126"
$ws.Range("L7").Value = [double]"35.83914203508134"
$ws.Range("M7").Value = "[12.397779235483895, 59.280504834678794]"
$ws.Range("N7").Value = [double]"0.003529976888952957"
$ws.Range("O7").Value = [double]"0.003529976888952957"
$ws.Range("P7").Value = [double]"-2.629000458859927"
$ws.Range("Q7").Value = "[-3.2956847857478517, -1.9623161319720022]"
$ws.Range("R7").Value = [double]"4.181401891401038e-10"
$ws.Range("S7").Value = [double]"4.181401891401038e-10"
$ws.Range("T7").Value = [double]"63.4764132591892"
$ws.Range("U7").Value = "[51.03812998788183, 75.91469653049657]"
$ws.Range("V7").Value = [double]"2.182698466413058e-13"
$ws.Range("W7").Value = [double]"2.182698466413058e-13"
$ws.Range("X7").Value = [double]"9.59015015015021"
$ws.Range("Y7").Value = [double]"7.158198198198241"
$ws.Range("Z7").Value = [double]"12.02210210210218"
